# Update header row (A1:B1) to use fixed, required column header names,
# and fix the "misery_index" header to use a space instead of an underscore.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "situation"
$ws.Range("B1").Value = "misery index"

# Move the active cell selection to B2 (was A4).
$ws.Range("B2").Select()

# Best-effort: reposition the workbook window on screen (cosmetic;
# xWindow/yWindow in the saved file are not affected by this runtime).
$wb.Windows.Item(1).Left = 1080
$wb.Windows.Item(1).Top = 1080
